# Commit: "Add learning records for 2025-03-03"
#
# The whole daily-log date column (A1:A366) shifts forward by 8 days
# (the log apparently wasn't touched for over a week, so the next entry
# lands on 2025-03-03 instead of 2025-02-23), and the newly-current day's
# row (row 2) gets its B/C notes filled in with "fawef" / "wef".
#
# All of column A stores plain text that merely looks like an ISO date
# (t="str" in the OOXML, no number format applied) - not a real Excel
# date serial. Writing a string like "2025-03-02" straight into .Value
# would get auto-coerced into a date by Excel, so each cell is forced to
# Text format before the write and then has that formatting override
# cleared again afterwards (ClearFormats doesn't revert the stored text
# back to a date - it only drops the formatting), which leaves the cells
# on the workbook's default style, matching the original file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowCount = 366
$startDate = Get-Date -Year 2025 -Month 3 -Day 2

$dataRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($rowCount, 1))

# Force text so the ISO-looking strings aren't auto-converted to dates.
$dataRange.NumberFormat = "@"

for ($r = 1; $r -le $rowCount; $r++) {
    $d = $startDate.AddDays($r - 1)
    $ws.Cells.Item($r, 1).Value = $d.ToString("yyyy-MM-dd")
}

# Drop the temporary text formatting so the cells stay on the default style.
$dataRange.ClearFormats()

# New entry for 2025-03-03 (row 2): learning-log note + mood/score text.
$ws.Cells.Item(2, 2).Value = "fawef"
$ws.Cells.Item(2, 3).Value = "wef"
